$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.52413133333333
$ws.Range("H2").Value = 37.572394
$ws.Range("I2").Value = 0.09718402715578596
$ws.Range("J2").Value = 0.1008592412859651
$ws.Range("M2").Value = 17.723347
$ws.Range("N2").Value = 53.170041
$ws.Range("O2").Value = 0.1083389314942055
$ws.Range("P2").Value = 0.1121884745845309
$ws.Range("Q2").Value = 221.9695254942393
$ws.Range("R2").Value = 1997.725729448154
$ws.Range("S2").Value = 0.0105288136603617
$ws.Range("T2").Value = 0.01131524442762557
$ws.Range("G3").Value = 12.52413133333333
$ws.Range("H3").Value = 37.572394
$ws.Range("I3").Value = 0.09718402715578596
$ws.Range("J3").Value = 0.1008592412859651
$ws.Range("O3").Value = 0.193467093096278
$ws.Range("P3").Value = 0.2003414447366573
$ws.Range("Q3").Value = 396.3838138428337
$ws.Range("R3").Value = 3567.454324585503
$ws.Range("S3").Value = 0.01880191122921965
$ws.Range("T3").Value = 0.02020628611427337
$ws.Range("G4").Value = 12.52413133333333
$ws.Range("H4").Value = 37.572394
$ws.Range("I4").Value = 0.09718402715578596
$ws.Range("J4").Value = 0.1008592412859651
$ws.Range("M4").Value = 47.45519633333333
$ws.Range("N4").Value = 142.365589
$ws.Range("O4").Value = 0.2900832029413559
$ws.Range("P4").Value = 0.3003905575931054
$ws.Range("Q4").Value = 594.3351113277851
$ws.Range("R4").Value = 5349.016001950065
$ws.Range("S4").Value = 0.0281914538720901
$ws.Range("T4").Value = 0.03029716372830862
$ws.Range("G5").Value = 12.52413133333333
$ws.Range("H5").Value = 37.572394
$ws.Range("I5").Value = 0.09718402715578596
$ws.Range("J5").Value = 0.1008592412859651
$ws.Range("M5").Value = 16.8400505
$ws.Range("N5").Value = 33.680101
$ws.Range("O5").Value = 0.102939533795646
$ws.Range("P5").Value = 0.07106481552351887
$ws.Range("Q5").Value = 210.9070041219656
$ws.Range("R5").Value = 1265.442024731794
$ws.Range("S5").Value = 0.0100040784478
$ws.Range("T5").Value = 0.00716754337582919
$ws.Range("G6").Value = 12.52413133333333
$ws.Range("H6").Value = 37.572394
$ws.Range("I6").Value = 0.09718402715578596
$ws.Range("J6").Value = 0.1008592412859651
$ws.Range("M6").Value = 49.92347333333333
$ws.Range("N6").Value = 149.77042
$ws.Range("O6").Value = 0.3051712386725145
$ws.Range("P6").Value = 0.3160147075621876
$ws.Range("Q6").Value = 625.248136642831
$ws.Range("R6").Value = 5627.233229785479
$ws.Range("S6").Value = 0.02965776994631449
$ws.Range("T6").Value = 0.03187300363992839
$ws.Range("I7").Value = 0.1842225641940495
$ws.Range("J7").Value = 0.1911893198517306
$ws.Range("M7").Value = 17.723347
$ws.Range("N7").Value = 53.170041
$ws.Range("O7").Value = 0.1083389314942055
$ws.Range("P7").Value = 0.1121884745845309
$ws.Range("Q7").Value = 420.7666255066346
$ws.Range("R7").Value = 3786.899629559712
$ws.Range("S7").Value = 0.019958475761906
$ws.Range("T7").Value = 0.02144923815101963
$ws.Range("I8").Value = 0.1842225641940495
$ws.Range("J8").Value = 0.1911893198517306
$ws.Range("O8").Value = 0.193467093096278
$ws.Range("P8").Value = 0.2003414447366573
$ws.Range("S8").Value = 0.03564100397736522
$ws.Range("T8").Value = 0.03830314455731459
$ws.Range("I9").Value = 0.1842225641940495
$ws.Range("J9").Value = 0.1911893198517306
$ws.Range("M9").Value = 47.45519633333333
$ws.Range("N9").Value = 142.365589
$ws.Range("O9").Value = 0.2900832029413559
$ws.Range("P9").Value = 0.3003905575931054
$ws.Range("Q9").Value = 1126.624831299161
$ws.Range("R9").Value = 10139.62348169245
$ws.Range("S9").Value = 0.05343987147547941
$ws.Range("T9").Value = 0.05743146639610794
$ws.Range("I10").Value = 0.1842225641940495
$ws.Range("J10").Value = 0.1911893198517306
$ws.Range("M10").Value = 16.8400505
$ws.Range("N10").Value = 33.680101
$ws.Range("O10").Value = 0.102939533795646
$ws.Range("P10").Value = 0.07106481552351887
$ws.Range("Q10").Value = 399.7964505376053
$ws.Range("R10").Value = 2398.778703225632
$ws.Range("S10").Value = 0.01896378487277392
$ws.Range("T10").Value = 0.01358683374533028
$ws.Range("I11").Value = 0.1842225641940495
$ws.Range("J11").Value = 0.1911893198517306
$ws.Range("M11").Value = 49.92347333333333
$ws.Range("N11").Value = 149.77042
$ws.Range("O11").Value = 0.3051712386725145
$ws.Range("P11").Value = 0.3160147075621876
$ws.Range("Q11").Value = 1185.223728229049
$ws.Range("R11").Value = 10667.01355406144
$ws.Range("S11").Value = 0.05621942810652489
$ws.Range("T11").Value = 0.0604186370019582
$ws.Range("G12").Value = 41.01852933333333
$ws.Range("H12").Value = 123.055588
$ws.Range("I12").Value = 0.3182932023406124
$ws.Range("J12").Value = 0.3303301152883236
$ws.Range("M12").Value = 17.723347
$ws.Range("N12").Value = 53.170041
$ws.Range("O12").Value = 0.1083389314942055
$ws.Range("P12").Value = 0.1121884745845309
$ws.Range("Q12").Value = 726.9856288043453
$ws.Range("R12").Value = 6542.870659239108
$ws.Range("S12").Value = 0.03448354544345088
$ws.Range("T12").Value = 0.03705923174352926
$ws.Range("G13").Value = 41.01852933333333
$ws.Range("H13").Value = 123.055588
$ws.Range("I13").Value = 0.3182932023406124
$ws.Range("J13").Value = 0.3303301152883236
$ws.Range("O13").Value = 0.193467093096278
$ws.Range("P13").Value = 0.2003414447366573
$ws.Range("Q13").Value = 1298.220264753756
$ws.Range("R13").Value = 11683.98238278381
$ws.Range("S13").Value = 0.06157926060914371
$ws.Range("T13").Value = 0.06617881253688931
$ws.Range("G14").Value = 41.01852933333333
$ws.Range("H14").Value = 123.055588
$ws.Range("I14").Value = 0.3182932023406124
$ws.Range("J14").Value = 0.3303301152883236
$ws.Range("M14").Value = 47.45519633333333
$ws.Range("N14").Value = 142.365589
$ws.Range("O14").Value = 0.2900832029413559
$ws.Range("P14").Value = 0.3003905575931054
$ws.Range("Q14").Value = 1946.542362817926
$ws.Range("R14").Value = 17518.88126536133
$ws.Range("S14").Value = 0.0923315116094259
$ws.Range("T14").Value = 0.0992280475212543
$ws.Range("G15").Value = 41.01852933333333
$ws.Range("H15").Value = 123.055588
$ws.Range("I15").Value = 0.3182932023406124
$ws.Range("J15").Value = 0.3303301152883236
$ws.Range("M15").Value = 16.8400505
$ws.Range("N15").Value = 33.680101
$ws.Range("O15").Value = 0.102939533795646
$ws.Range("P15").Value = 0.07106481552351887
$ws.Range("Q15").Value = 690.7541054090647
$ws.Range("R15").Value = 4144.524632454388
$ws.Range("S15").Value = 0.03276495385926585
$ws.Range("T15").Value = 0.02347484870482744
$ws.Range("G16").Value = 41.01852933333333
$ws.Range("H16").Value = 123.055588
$ws.Range("I16").Value = 0.3182932023406124
$ws.Range("J16").Value = 0.3303301152883236
$ws.Range("M16").Value = 49.92347333333333
$ws.Range("N16").Value = 149.77042
$ws.Range("O16").Value = 0.3051712386725145
$ws.Range("P16").Value = 0.3160147075621876
$ws.Range("Q16").Value = 2047.787455345218
$ws.Range("R16").Value = 18430.08709810696
$ws.Range("S16").Value = 0.09713393081932596
$ws.Range("T16").Value = 0.1043891747818233
$ws.Range("G17").Value = 14.087727
$ws.Range("H17").Value = 28.175454
$ws.Range("I17").Value = 0.1093171260259301
$ws.Range("J17").Value = 0.07563411885139956
$ws.Range("M17").Value = 17.723347
$ws.Range("N17").Value = 53.170041
$ws.Range("O17").Value = 0.1083389314942055
$ws.Range("P17").Value = 0.1121884745845309
$ws.Range("Q17").Value = 249.681674062269
$ws.Range("R17").Value = 1498.090044373614
$ws.Range("S17").Value = 0.01184330062766666
$ws.Range("T17").Value = 0.00848527642048363
$ws.Range("G18").Value = 14.087727
$ws.Range("H18").Value = 28.175454
$ws.Range("I18").Value = 0.1093171260259301
$ws.Range("J18").Value = 0.07563411885139956
$ws.Range("O18").Value = 0.193467093096278
$ws.Range("P18").Value = 0.2003414447366573
$ws.Range("Q18").Value = 445.870999593744
$ws.Range("R18").Value = 2675.225997562464
$ws.Range("S18").Value = 0.02114926659787617
$ws.Range("T18").Value = 0.01515264864207344
$ws.Range("G19").Value = 14.087727
$ws.Range("H19").Value = 28.175454
$ws.Range("I19").Value = 0.1093171260259301
$ws.Range("J19").Value = 0.07563411885139956
$ws.Range("M19").Value = 47.45519633333333
$ws.Range("N19").Value = 142.365589
$ws.Range("O19").Value = 0.2900832029413559
$ws.Range("P19").Value = 0.3003905575931054
$ws.Range("Q19").Value = 668.5358506754011
$ws.Range("R19").Value = 4011.215104052406
$ws.Range("S19").Value = 0.03171106205394564
$ws.Range("T19").Value = 0.02271977513483512
$ws.Range("G20").Value = 14.087727
$ws.Range("H20").Value = 28.175454
$ws.Range("I20").Value = 0.1093171260259301
$ws.Range("J20").Value = 0.07563411885139956
$ws.Range("M20").Value = 16.8400505
$ws.Range("N20").Value = 33.680101
$ws.Range("O20").Value = 0.102939533795646
$ws.Range("P20").Value = 0.07106481552351887
$ws.Range("Q20").Value = 237.2380341102135
$ws.Range("R20").Value = 948.9521364408541
$ws.Range("S20").Value = 0.01125305398898912
$ws.Range("T20").Value = 0.00537492470345861
$ws.Range("G21").Value = 14.087727
$ws.Range("H21").Value = 28.175454
$ws.Range("I21").Value = 0.1093171260259301
$ws.Range("J21").Value = 0.07563411885139956
$ws.Range("M21").Value = 49.92347333333333
$ws.Range("N21").Value = 149.77042
$ws.Range("O21").Value = 0.3051712386725145
$ws.Range("P21").Value = 0.3160147075621876
$ws.Range("Q21").Value = 703.3082632117801
$ws.Range("R21").Value = 4219.84957927068
$ws.Range("S21").Value = 0.03336044275745245
$ws.Range("T21").Value = 0.02390149395054877
$ws.Range("G22").Value = 37.49906666666667
$ws.Range("H22").Value = 112.4972
$ws.Range("I22").Value = 0.2909830802836222
$ws.Range("J22").Value = 0.3019872047225811
$ws.Range("M22").Value = 17.723347
$ws.Range("N22").Value = 53.170041
$ws.Range("O22").Value = 0.1083389314942055
$ws.Range("P22").Value = 0.1121884745845309
$ws.Range("Q22").Value = 664.6089707094668
$ws.Range("R22").Value = 5981.4807363852
$ws.Range("S22").Value = 0.03152479600082024
$ws.Range("T22").Value = 0.03387948384187283
$ws.Range("G23").Value = 37.49906666666667
$ws.Range("H23").Value = 112.4972
$ws.Range("I23").Value = 0.2909830802836222
$ws.Range("J23").Value = 0.3019872047225811
$ws.Range("O23").Value = 0.193467093096278
$ws.Range("P23").Value = 0.2003414447366573
$ws.Range("Q23").Value = 1186.830660368356
$ws.Range("R23").Value = 10681.4759433152
$ws.Range("S23").Value = 0.05629565068267329
$ws.Range("T23").Value = 0.0605005528861066
$ws.Range("G24").Value = 37.49906666666667
$ws.Range("H24").Value = 112.4972
$ws.Range("I24").Value = 0.2909830802836222
$ws.Range("J24").Value = 0.3019872047225811
$ws.Range("M24").Value = 47.45519633333333
$ws.Range("N24").Value = 142.365589
$ws.Range("O24").Value = 0.2900832029413559
$ws.Range("P24").Value = 0.3003905575931054
$ws.Range("Q24").Value = 1779.525570983423
$ws.Range("R24").Value = 16015.7301388508
$ws.Range("S24").Value = 0.08440930393041483
$ws.Range("T24").Value = 0.09071410481259941
$ws.Range("G25").Value = 37.49906666666667
$ws.Range("H25").Value = 112.4972
$ws.Range("I25").Value = 0.2909830802836222
$ws.Range("J25").Value = 0.3019872047225811
$ws.Range("M25").Value = 16.8400505
$ws.Range("N25").Value = 33.680101
$ws.Range("O25").Value = 0.102939533795646
$ws.Range("P25").Value = 0.07106481552351887
$ws.Range("Q25").Value = 631.4861763695334
$ws.Range("R25").Value = 3788.9170582172
$ws.Range("S25").Value = 0.0299536626268171
$ws.Range("T25").Value = 0.02146066499407335
$ws.Range("G26").Value = 37.49906666666667
$ws.Range("H26").Value = 112.4972
$ws.Range("I26").Value = 0.2909830802836222
$ws.Range("J26").Value = 0.3019872047225811
$ws.Range("M26").Value = 49.92347333333333
$ws.Range("N26").Value = 149.77042
$ws.Range("O26").Value = 0.3051712386725145
$ws.Range("P26").Value = 0.3160147075621876
$ws.Range("Q26").Value = 1872.083654758223
$ws.Range("R26").Value = 16848.752892824
$ws.Range("S26").Value = 0.08879966704289673
$ws.Range("T26").Value = 0.09543239818792894
